# Magic Beans Demo Presentation — "Work on fixing runtime errors"
#
# 1. Bump the cached "today" date field (datetimeFigureOut) on the slide
#    master and every slide layout from 12/14/2015 -> 12/15/2015.
# 2. Append a new slide 7 ("Issues") after the existing 6 slides, using
#    the "Title and Content" layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date field text wherever it appears.
# ---------------------------------------------------------------------
$newDate = "12/15/2015"

$master = $p.SlideMaster

foreach ($sh in $master.Shapes) {
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    foreach ($sh in $layout.Shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Add the new "Issues" slide at the end (index 7), Title+Content
#    layout (layout 2 == ppLayoutText / "Title and Content").
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Issues"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Understanding "
[void]$body.InsertAfter("LaTeX")
[void]$body.InsertAfter(", what it is, how it works")
[void]$body.InsertAfter("`rLearning the ")
[void]$body.InsertAfter("Netbeans")
[void]$body.InsertAfter(" architecture proved very hard")
[void]$body.InsertAfter("`rBecause of this choosing an approach was difficult and we had to change directions 3+ times which came with re-learning")
[void]$body.InsertAfter("`rConnecting our plugin to ")
[void]$body.InsertAfter("LaTeX")
[void]$body.InsertAfter("`rGetting the rendered PDF back from the executable")
[void]$body.InsertAfter("`r")

# "Getting the rendered PDF back from the executable" is a sub-bullet
# (second outline level).
$body.Paragraphs(5).IndentLevel = 2

Write-Output "done"
